$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values that look numeric are written as text,
# matching the source workbook where column D cells are stored as inline strings.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.021.92"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.553.65"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.28"
$ws.Range("E5").Value = "  +1.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.35"
$ws.Range("E6").Value = "  -1.95%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.97%  "

# Row 9
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.52"
$ws.Range("E10").Value = "  -4.24%  "

# Row 11
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("E12").Value = "  -1.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.21"
$ws.Range("E13").Value = "  -3.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.010.26"
$ws.Range("E14").Value = "  -0.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.932.03"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("E16").Value = "  -1.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.549.57"
$ws.Range("E17").Value = "  -0.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.34"
$ws.Range("E18").Value = "  -3.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "336.47"
$ws.Range("E19").Value = "  -1.74%  "

# Row 20
$ws.Range("E20").Value = "  -1.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.54"
$ws.Range("E23").Value = "  -0.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.62"
$ws.Range("E25").Value = "  +2.06%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.48"
$ws.Range("E27").Value = "  +0.83%  "

# Row 28
$ws.Range("E28").Value = "  -1.21%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  +1.84%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  +1.53%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0813"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "177.45"
$ws.Range("E32").Value = "  -0.18%  "

# Row 33
$ws.Range("E33").Value = "  -1.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "411.91"
$ws.Range("E34").Value = "  -0.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.14"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("E36").Value = "  -2.11%  "

# Row 37
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("E38").Value = "  -2.70%  "

# Row 39
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.88"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.08"
$ws.Range("E42").Value = "  -2.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  -1.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.85"
$ws.Range("E44").Value = "  -1.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0537"
$ws.Range("E45").Value = "  +0.54%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  -1.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0968"

# Row 48
$ws.Range("E48").Value = "  +1.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.26"
$ws.Range("E49").Value = "  -2.91%  "

# Row 50
$ws.Range("E50").Value = "  -7.94%  "

# Row 51
$ws.Range("E51").Value = "  -0.08%  "
